$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Objetivos: value replaced with the docente's name/code
$ws.Range("B10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

# Row 13: becomes "Programa resumido:" / "Semestral" (with 60pt row height)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# Row 14: "Short syllabus:" row stays the same text, keep as-is (60pt), but re-assert for safety
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Introduction. Soil formation. Soil physical properties. Soil classification. Soil water. Practical class: Profile description in the field. Practical class: Characterization and determination methods of physical and hydraulic properties of the soil."
$ws.Range("C14").Value = "Introduction. Soil formation. Soil physical properties. Soil classification. Soil water. Practical class: Profile description in the field. Practical class: Characterization and determination methods of physical and hydraulic properties of the soil."
$ws.Rows(14).RowHeight = 60

# Row 15: becomes "Programa:" / "01/01/2020" (120pt row height)
# (copy the already-text-typed "01/01/2020" cell from row 8 so it is not
# reinterpreted as a date serial number)
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))
$ws.Rows(15).RowHeight = 120

# Row 16: "Syllabus:" row - unchanged content, re-assert
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "INTRODUCTION. Basic Concepts. The soil profile. Horizons and layers definition and notation. SOIL FORMATION. Formation and factors processes. Weathering. PHYSICAL ATTRIBUTES OF THE SOIL. Volumetric composition, grain size and texture, structure and aggregation, color, porosity, density and compression, consistency. SOIL CLASSIFICATION. Brazilian system of soil classification. Main morphological attributes. Top Soil Classes. SOIL WATER. Concept and importance. Moisture constants. Total potential of Soil water and its components. Characteristic curve of soil water. Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil."
$ws.Range("C16").Value = "INTRODUCTION. Basic Concepts. The soil profile. Horizons and layers definition and notation. SOIL FORMATION. Formation and factors processes. Weathering. PHYSICAL ATTRIBUTES OF THE SOIL. Volumetric composition, grain size and texture, structure and aggregation, color, porosity, density and compression, consistency. SOIL CLASSIFICATION. Brazilian system of soil classification. Main morphological attributes. Top Soil Classes. SOIL WATER. Concept and importance. Moisture constants. Total potential of Soil water and its components. Characteristic curve of soil water. Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil."
$ws.Rows(16).RowHeight = 120

# Row 17: becomes "Avaliação:" only (A column), no B/C -> default row height
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Delete() | Out-Null
$ws.Rows(17).AutoFit() | Out-Null

# Row 18: becomes "Método:" / docente text (60pt row height)
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Rows(18).RowHeight = 60

# Row 19: becomes "Critério:" / avaliação description (60pt row height)
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Rows(19).RowHeight = 60

# Row 20: becomes "Norma de recuperação:" / critérios description (60pt row height)
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Rows(20).RowHeight = 60

# Row 21: becomes "Bibliografia:" / exame final description (120pt row height)
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
$ws.Range("C21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
$ws.Rows(21).RowHeight = 120

# Row 22 (old "Bibliografia:" / "Bibliografia básica..." row) is removed entirely.
$ws.Rows(22).Delete() | Out-Null
